$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells P1/Q1 should match the style of the existing header row
# (bold font + border), so copy formatting from O1 before writing the values.
$ws.Range("O1").Copy($ws.Range("P1"))
$ws.Range("O1").Copy($ws.Range("Q1"))
$ws.Range("P1").Value = 14
$ws.Range("Q1").Value = 15

# For rows 2 through 25: swap I<->K values, swap M<->O values, and add P/Q = 2
for ($r = 2; $r -le 25; $r++) {
    $iVal = $ws.Cells.Item($r, 9).Value2   # I
    $kVal = $ws.Cells.Item($r, 11).Value2  # K
    $mVal = $ws.Cells.Item($r, 13).Value2  # M
    $oVal = $ws.Cells.Item($r, 15).Value2  # O

    $ws.Cells.Item($r, 9).Value = $kVal    # I = old K
    $ws.Cells.Item($r, 11).Value = $iVal   # K = old I
    $ws.Cells.Item($r, 13).Value = $oVal   # M = old O
    $ws.Cells.Item($r, 15).Value = $mVal   # O = old M

    $ws.Cells.Item($r, 16).Value = 2       # P
    $ws.Cells.Item($r, 17).Value = 2       # Q
}
